# Generate Report for Handoff
# A new handoff was generated for the "b6dece24-..." file (row 5 in each sheet),
# so its "Latest Handoff Date" / "Latest Handoff Datetime" values are refreshed.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-03-23 06:04:05"

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-23 06:03:56"

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-23 06:04:05"
